# Remove the "Modified" info row from the "ApPredict version information" sheet.
$wb = $excel.ActiveWorkbook
$originalActive = $wb.ActiveSheet
$ws = $wb.Worksheets.Item("ApPredict version information")

# Delete the entire row 3 (the "Modified" / TRUE row), shifting rows below up.
$ws.Rows.Item(3).Delete()

# Select the new row 3 (now "Build options"), matching post-edit selection state.
$ws.Activate()
$ws.Range("A3:XFD3").Select()

# Restore the workbook's originally active sheet so the overall active tab is unchanged.
$originalActive.Activate()
